$d = $word.ActiveDocument

# --- Fix "Ksiegarnia" typo -> "Księgarnia" -----------------------------
$d.Content.Find.Execute("Ksiegarnia", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Księgarnia", 2) | Out-Null

# --- Merge the name runs that were split for spell-check purposes ------
# (no textual change, just collapses "Artur "+"Pereć" / "Mateusz "+"Malisz"
#  into single contiguous runs, matching how the doc now reads)
$d.Content.Find.Execute("Artur Pereć", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Artur Pereć", 2) | Out-Null
$d.Content.Find.Execute("Mateusz Malisz", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Mateusz Malisz", 2) | Out-Null

# --- Artur Pereć: add his GitHub handle and the review date ------------
$r = $d.Paragraphs(5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" (artur1947)")
$r.Collapse(0)
$r.InsertAfter(" – 22.11.2016r.")

# --- Mateusz Malisz: add his GitHub handle and the review date ---------
$r = $d.Paragraphs(6).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" (MatthewRock)")
$r.Collapse(0)
$r.InsertAfter(" – 10.01.2017r.")

# --- Piotr Radwan: add his GitHub handle and the review date -----------
$r = $d.Paragraphs(7).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" (Antystenes)")
$r.Collapse(0)
$r.InsertAfter(" – 06.12.2016r.")

# --- Patrycja Stefańska: add her GitHub handle and the review date -----
$r = $d.Paragraphs(8).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" (malaczarna)")
$r.Collapse(0)
$r.InsertAfter(" – 13.12.2016r.")
